# Update cryptos table: latest Price (col D) and Volume(1h) (col E) per source refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.091.75"
$ws.Range("E2").Value = "  +0.49%  "
$ws.Range("D3").Value = "2.328.48"
$ws.Range("E3").Value = "  +4.06%  "
$ws.Range("E4").Value = "  -0.03%  "
$c = $ws.Range("D5")
$c.Value = "'98.61"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +4.71%  "
$c = $ws.Range("D6")
$c.Value = "'271.98"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +0.70%  "
$c = $ws.Range("D7")
$c.Value = "'0.631"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +0.35%  "
$ws.Range("E8").Value = "  -0.06%  "
$c = $ws.Range("D9")
$c.Value = "'0.627"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -1.68%  "
$c = $ws.Range("D10")
$c.Value = "'45.84"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -0.61%  "
$c = $ws.Range("D11")
$c.Value = "'0.0957"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -0.10%  "
$ws.Range("E12").Value = "  -4.79%  "
$ws.Range("E13").Value = "  +0.46%  "
$ws.Range("D14").Value = "2.668.37"
$ws.Range("E14").Value = "  +3.59%  "
$c = $ws.Range("D15")
$c.Value = "'15.56"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +1.71%  "
$c = $ws.Range("D16")
$c.Value = "'0.879"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +7.24%  "
$ws.Range("D17").Value = "2.332.49"
$ws.Range("E17").Value = "  +4.70%  "
$ws.Range("D18").Value = "44.036.09"
$ws.Range("E18").Value = "  +0.34%  "
$ws.Range("E19").Value = "  +5.03%  "
$c = $ws.Range("D20")
$c.Value = "'6.43"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +3.75%  "
$c = $ws.Range("D21")
$c.Value = "'73.68"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +4.05%  "
$ws.Range("E22").Value = "  -1.19%  "
$c = $ws.Range("D23")
$c.Value = "'240.49"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +2.48%  "
$c = $ws.Range("D24")
$c.Value = "'9.32"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +2.34%  "
$ws.Range("E25").Value = "  -0.05%  "
$ws.Range("E26").Value = "  +1.81%  "
$c = $ws.Range("D27")
$c.Value = "'11.45"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +0.39%  "
$c = $ws.Range("D28")
$c.Value = "'3.51"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -1.94%  "
$c = $ws.Range("D29")
$c.Value = "'2.31"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +2.22%  "
$c = $ws.Range("D30")
$c.Value = "'38.35"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -5.18%  "
$c = $ws.Range("D31")
$c.Value = "'22.46"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +6.85%  "
$c = $ws.Range("D32")
$c.Value = "'175.76"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +1.75%  "
$c = $ws.Range("D33")
$c.Value = "'0.0918"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +0.46%  "
$c = $ws.Range("D34")
$c.Value = "'5.53"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +2.12%  "
$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$c = $ws.Range("D36")
$c.Value = "'0.110"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -2.07%  "
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$c = $ws.Range("D37")
$c.Value = "'0.0364"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +3.33%  "
$ws.Range("E38").Value = "  +4.13%  "
$ws.Range("E39").Value = "  -4.84%  "
$ws.Range("E40").Value = "  +8.90%  "
$c = $ws.Range("D41")
$c.Value = "'2.42"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +12.12%  "
$c = $ws.Range("D42")
$c.Value = "'1.44"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +24.30%  "
$ws.Range("E43").Value = "  -2.61%  "
$c = $ws.Range("D44")
$c.Value = "'62.93"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -0.70%  "
$ws.Range("E45").Value = "  +8.67%  "
$c = $ws.Range("D46")
$c.Value = "'5.36"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -1.04%  "
$ws.Range("E47").Value = "  +4.37%  "
$c = $ws.Range("D48")
$c.Value = "'100.59"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -0.76%  "
$ws.Range("E49").Value = "  +0.90%  "
$ws.Range("E50").Value = "  +16.76%  "
$ws.Range("D51").Value = "2.554.92"
$ws.Range("E51").Value = "  +3.93%  "
